$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G18 was stored as text "26.0" - convert it to the real number 26
$ws.Range("G18").Value = 26

# Append the new row 19 of pool-log data
$ws.Range("A19").Value = "'2026-02-20"
$ws.Range("A19").Style = "Normal"

$ws.Range("B19").Value = "Sol"
$ws.Range("C19").Value = "Agradavel"
$ws.Range("D19").Value = "normal"
$ws.Range("E19").Value = "aula"
$ws.Range("F19").Value = "nenhuma"

$ws.Range("G19").Value = "'27"
$ws.Range("G19").Style = "Normal"

$ws.Range("H19").Value = "'32"
$ws.Range("H19").Style = "Normal"

$ws.Range("I19").Value = 2.5

$ws.Range("J19").Value = "dqs01"
$ws.Range("K19").Value = "Quarta e Sexta"
$ws.Range("L19").Value = "08:00"
$ws.Range("M19").Value = "Daniela"
